# Apply the change described by the diff:
#   Insert a new "culture_collection" attribute column before column U
#   on the header row (row 15), shifting U:BV one column to the right
#   (to V:BW), and add its header comment. All pre-existing comments
#   must follow their cells to the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColToNum($colLetters) {
    $numResult = 0
    for ($ctnIdx = 0; $ctnIdx -lt $colLetters.Length; $ctnIdx++) {
        $ctnCh = $colLetters.Substring($ctnIdx, 1)
        $numResult = $numResult * 26 + ([int][char]$ctnCh - [int][char]'A' + 1)
    }
    return $numResult
}

function NumToCol($colNumber) {
    $strResult = ""
    $ntcRemaining = $colNumber
    while ($ntcRemaining -gt 0) {
        $ntcRem = ($ntcRemaining - 1) % 26
        $strResult = [char]([int][char]'A' + $ntcRem) + $strResult
        $ntcRemaining = [int](($ntcRemaining - 1) / 26)
    }
    return $strResult
}

$insertColNum = ColToNum "U"

# ---- Step 1: capture every existing comment (address + text) up front,
#      before anything on the sheet is touched. ----
$commentCount = $ws.Comments.Count()
$capturedAddr = @()
$capturedText = @()
for ($captureIdx = 1; $captureIdx -le $commentCount; $captureIdx++) {
    $existingComment = $ws.Comments.Item($captureIdx)
    $capturedAddr += $existingComment.Parent().Address()
    $capturedText += $existingComment.Text()
}

# ---- Step 2: remove all existing comments (we'll re-create them at the
#      (possibly shifted) correct location afterwards). ----
for ($deleteIdx = $commentCount; $deleteIdx -ge 1; $deleteIdx--) {
    $ws.Comments.Item($deleteIdx).Delete()
}

# ---- Step 3: insert the new column at U, shifting U:BV -> V:BW. ----
$ws.Columns("U:U").Insert()

# ---- Step 4: give the new column its header text. ----
$ws.Range("U15").Value = "culture_collection"

# ---- Step 5: re-create the captured comments, shifting any comment that
#      was in column U or later one column to the right. ----
for ($rebuildIdx = 0; $rebuildIdx -lt $capturedAddr.Length; $rebuildIdx++) {
    $rawAddr = $capturedAddr[$rebuildIdx]
    $cleanAddr = $rawAddr.Replace('$', '')
    if ($cleanAddr -match '^([A-Z]+)([0-9]+)$') {
        $oldColLetters = $matches[1]
        $rowNumber = $matches[2]
        $oldColNum = ColToNum $oldColLetters
        $newColNum = $oldColNum
        if ($oldColNum -ge $insertColNum) {
            $newColNum = $oldColNum + 1
        }
        $newColLetters = NumToCol $newColNum
        $newAddr = $newColLetters + $rowNumber
        $ws.Range($newAddr).AddComment($capturedText[$rebuildIdx])
    }
}

# ---- Step 6: add the brand-new comment for the culture_collection header. ----
$ws.Range("U15").AddComment("Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier")
